$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B2" = 0.1893939393939394
    "C2" = 0.5568181818181818
    "J2" = 0.007575757575757576
    "P2" = 0.1022727272727273
    "S2" = 0.143939393939394
    "C3" = 0.006666666666666667
    "J3" = 0.02
    "P3" = 0.7733333333333333
    "S3" = 0.2
    "P4" = 0.6842105263157895
    "S4" = 0.3157894736842105
    "B6" = 0.1056910569105691
    "D6" = 0.004065040650406504
    "F6" = 0.06097560975609756
    "J6" = 0.2357723577235772
    "O6" = 0.01626016260162602
    "Q6" = 0.1747967479674797
    "R6" = 0.06910569105691057
    "S6" = 0.3333333333333333
    "B7" = 0.04147465437788019
    "D7" = 0.03225806451612903
    "F7" = 0.06451612903225806
    "J7" = 0.1428571428571428
    "O7" = 0.009216589861751152
    "Q7" = 0.2304147465437788
    "R7" = 0.05069124423963134
    "S7" = 0.4285714285714285
    "B8" = 0.08281573498964803
    "D8" = 0.01656314699792961
    "F8" = 0.06625258799171843
    "J8" = 0.08074534161490683
    "O8" = 0.02070393374741201
    "Q8" = 0.1966873706004141
    "R8" = 0.09316770186335403
    "S8" = 0.443064182194617
    "B9" = 0.07623318385650224
    "D9" = 0.01345291479820628
    "F9" = 0.04932735426008968
    "J9" = 0.09417040358744394
    "O9" = 0.0179372197309417
    "Q9" = 0.1883408071748879
    "R9" = 0.1031390134529148
    "S9" = 0.4573991031390134
    "B10" = 0.09754098360655737
    "D10" = 0.02131147540983606
    "E10" = 0.000819672131147541
    "F10" = 0.0860655737704918
    "J10" = 0.08852459016393442
    "O10" = 0.02295081967213115
    "Q10" = 0.2180327868852459
    "R10" = 0.08278688524590164
    "S10" = 0.3819672131147541
    "G11" = 0.09810126582278481
    "J11" = 0.1044303797468354
    "K11" = 0.1645569620253164
    "L11" = 0.620253164556962
    "S11" = 0.01265822784810127
    "G12" = 0.7884615384615384
    "J12" = 0.1634615384615385
    "L12" = 0.02884615384615385
    "S12" = 0.01923076923076923
    "G13" = 0.6111111111111112
    "J13" = 0.3148148148148148
    "S13" = 0.07407407407407407
    "F15" = 0.02024291497975709
    "H15" = 0.1862348178137652
    "I15" = 0.0728744939271255
    "J15" = 0.319838056680162
    "K15" = 0.08502024291497975
    "O15" = 0.05668016194331984
    "S15" = 0.2591093117408907
    "H16" = 0.1666666666666667
    "I16" = 0.07407407407407407
    "J16" = 0.4444444444444444
    "K16" = 0.1172839506172839
    "M16" = 0.0308641975308642
    "O16" = 0.06172839506172839
    "S16" = 0.1049382716049383
    "F17" = 0.02469135802469136
    "H17" = 0.1707818930041152
    "I17" = 0.102880658436214
    "J17" = 0.411522633744856
    "K17" = 0.102880658436214
    "M17" = 0.02880658436213992
    "O17" = 0.07407407407407407
    "S17" = 0.08436213991769548
    "F18" = 0.01020408163265306
    "H18" = 0.2091836734693878
    "I18" = 0.09693877551020408
    "J18" = 0.3877551020408163
    "K18" = 0.1020408163265306
    "M18" = 0.02040816326530612
    "O18" = 0.06122448979591837
    "S18" = 0.1122448979591837
    "F19" = 0.01948558067030397
    "H19" = 0.2244738893219018
    "I19" = 0.09586905689789556
    "J19" = 0.3616523772408418
    "K19" = 0.1153546375681995
    "M19" = 0.02494154325798909
    "O19" = 0.0740452065471551
    "S19" = 0.08417770849571317
}

foreach ($cellRef in $changes.Keys) {
    $ws.Range($cellRef).Value = $changes[$cellRef]
}

Write-Output "Updated $($changes.Count) cells"